$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.465.22"
$ws.Range("E2").Value = "  -1.15%  "
$ws.Range("D3").Value = "2.915.94"
$ws.Range("E3").Value = "  +1.67%  "
$ws.Range("E4").Value = "  +0.11%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "352.32"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +0.43%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "109.69"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -2.13%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.557"
$c.Style = "Normal"
$ws.Range("E7").Value = "  -0.41%  "
$ws.Range("E8").Value = "  +0.00%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.623"
$c.Style = "Normal"
$ws.Range("E9").Value = "  +0.42%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "38.69"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -3.72%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.0895"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +5.05%  "
$ws.Range("E12").Value = "  +0.62%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "19.56"
$c.Style = "Normal"
$ws.Range("E13").Value = "  -2.47%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "7.92"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +1.43%  "
$ws.Range("D15").Value = "3.379.70"
$ws.Range("E15").Value = "  +1.84%  "
$ws.Range("D16").Value = "2.913.56"
$ws.Range("E16").Value = "  +2.02%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "0.979"
$c.Style = "Normal"
$ws.Range("E17").Value = "  -1.26%  "
$ws.Range("D18").Value = "51.538.78"
$ws.Range("E18").Value = "  -0.98%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "7.54"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -1.45%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "14.16"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +4.07%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "3.22"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -4.03%  "
$ws.Range("D22").Value = "0.0₃0979"
$ws.Range("E22").Value = "  +0.14%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "70.73"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -0.09%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "269.31"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +0.05%  "
$ws.Range("E25").Value = "  -0.03%  "
$ws.Range("E26").Value = "  +8.41%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "27.08"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +2.94%  "
$ws.Range("E28").Value = "  +0.31%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "7.44"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +18.92%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "0.108"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +17.09%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "10.56"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -0.60%  "
$ws.Range("B32").Value = "RenderToken"
$ws.Range("C32").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "6.16"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +3.10%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "36.46"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -6.23%  "
$ws.Range("B34").Value = "OKB"
$ws.Range("C34").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "52.11"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -1.13%  "
$ws.Range("B35").Value = "VeChain"
$ws.Range("C35").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "0.0435"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -5.13%  "
$ws.Range("B36").Value = "FirstDigitalUSD"
$ws.Range("C36").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -0.07%  "
$ws.Range("B37").Value = "Toncoin"
$ws.Range("C37").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "1.92"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -15.11%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "3.23"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -1.96%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "18.20"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -1.82%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "2.01"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -0.73%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "2.66"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +2.72%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "0.117"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +0.52%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "23.31"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +4.38%  "
$ws.Range("E44").Value = "  -1.25%  "
$ws.Range("E45").Value = "  +1.67%  "
$ws.Range("B46").Value = "Maker"
$ws.Range("C46").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D46").Value = "2.136.45"
$ws.Range("E46").Value = "  -1.93%  "
$ws.Range("B47").Value = "Monero"
$ws.Range("C47").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "113.51"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -6.50%  "
$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "3.43"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -3.96%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "0.247"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +0.05%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "0.0326"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +1.77%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "9.03"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +0.07%  "
